$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NBA games added for February 12, 2025 (rows 800-814)

# Row 800: ORL @ CHA
$ws.Cells.Item(800, 1).Value = 45700
$ws.Cells.Item(800, 2).Value = "Wednesday"
$ws.Cells.Item(800, 3).Value = "CHA"
$ws.Cells.Item(800, 4).Value = "ORL"
$ws.Cells.Item(800, 5).Value = 0
$ws.Cells.Item(800, 6).Value = "Curtis Blair"
$ws.Cells.Item(800, 7).Value = "Ray Acosta"
$ws.Cells.Item(800, 8).Value = "Robert Hussey"
$ws.Cells.Item(800, 9).Value = 208
$ws.Cells.Item(800, 10).Value = -10.5
$ws.Cells.Item(800, 11).Value = 86
$ws.Cells.Item(800, 12).Value = 102
$ws.Cells.Item(800, 13).Formula = "=K800+L800"
$ws.Cells.Item(800, 14).Formula = "=(L800-K800)*-1"
$ws.Cells.Item(800, 15).Value = 1
$ws.Cells.Item(800, 16).Formula = "=IF(M800>I800,1,0)"
$ws.Cells.Item(800, 17).Formula = "=IF(P800=1,(M800-I800), """")"
$ws.Cells.Item(800, 18).Formula = "=IF(M800<I800, 1, 0)"
$ws.Cells.Item(800, 19).Formula = "=IF(R800=1,(I800-M800),"""")"
$ws.Cells.Item(800, 20).Formula = "=IF(M800=I800,1,0)"

# Row 801: WAS @ IND
$ws.Cells.Item(801, 1).Value = 45700
$ws.Cells.Item(801, 2).Value = "Wednesday"
$ws.Cells.Item(801, 3).Value = "IND"
$ws.Cells.Item(801, 4).Value = "WAS"
$ws.Cells.Item(801, 5).Value = 0
$ws.Cells.Item(801, 6).Value = "John Goble"
$ws.Cells.Item(801, 7).Value = "Nick Buchert"
$ws.Cells.Item(801, 8).Value = "Danielle Scott"
$ws.Cells.Item(801, 9).Value = 239.5
$ws.Cells.Item(801, 10).Value = 10
$ws.Cells.Item(801, 11).Value = 134
$ws.Cells.Item(801, 12).Value = 130
$ws.Cells.Item(801, 13).Formula = "=K801+L801"
$ws.Cells.Item(801, 14).Formula = "=(L801-K801)*-1"
$ws.Cells.Item(801, 15).Value = 1
$ws.Cells.Item(801, 16).Formula = "=IF(M801>I801,1,0)"
$ws.Cells.Item(801, 17).Formula = "=IF(P801=1,(M801-I801), """")"
$ws.Cells.Item(801, 18).Formula = "=IF(M801<I801, 1, 0)"
$ws.Cells.Item(801, 19).Formula = "=IF(R801=1,(I801-M801),"""")"
$ws.Cells.Item(801, 20).Formula = "=IF(M801=I801,1,0)"

# Row 802: BOS @ SAS
$ws.Cells.Item(802, 1).Value = 45700
$ws.Cells.Item(802, 2).Value = "Wednesday"
$ws.Cells.Item(802, 3).Value = "SAS"
$ws.Cells.Item(802, 4).Value = "BOS"
$ws.Cells.Item(802, 5).Value = 0
$ws.Cells.Item(802, 6).Value = "Kevin Scott"
$ws.Cells.Item(802, 7).Value = "Phenizee Ransom"
$ws.Cells.Item(802, 8).Value = "Derrick Collins"
$ws.Cells.Item(802, 9).Value = 231.5
$ws.Cells.Item(802, 10).Value = -9.5
$ws.Cells.Item(802, 11).Value = 103
$ws.Cells.Item(802, 12).Value = 116
$ws.Cells.Item(802, 13).Formula = "=K802+L802"
$ws.Cells.Item(802, 14).Formula = "=(L802-K802)*-1"
$ws.Cells.Item(802, 15).Value = 1
$ws.Cells.Item(802, 16).Formula = "=IF(M802>I802,1,0)"
$ws.Cells.Item(802, 17).Formula = "=IF(P802=1,(M802-I802), """")"
$ws.Cells.Item(802, 18).Formula = "=IF(M802<I802, 1, 0)"
$ws.Cells.Item(802, 19).Formula = "=IF(R802=1,(I802-M802),"""")"
$ws.Cells.Item(802, 20).Formula = "=IF(M802=I802,1,0)"

# Row 803: BKN @ PHI
$ws.Cells.Item(803, 1).Value = 45700
$ws.Cells.Item(803, 2).Value = "Wednesday"
$ws.Cells.Item(803, 3).Value = "PHI"
$ws.Cells.Item(803, 4).Value = "BKN"
$ws.Cells.Item(803, 5).Value = 0
$ws.Cells.Item(803, 6).Value = "Brian Forte"
$ws.Cells.Item(803, 7).Value = "Eric Dalen"
$ws.Cells.Item(803, 8).Value = "Suyash Mehta"
$ws.Cells.Item(803, 9).Value = 212
$ws.Cells.Item(803, 10).Value = 2
$ws.Cells.Item(803, 11).Value = 96
$ws.Cells.Item(803, 12).Value = 100
$ws.Cells.Item(803, 13).Formula = "=K803+L803"
$ws.Cells.Item(803, 14).Formula = "=(L803-K803)*-1"
$ws.Cells.Item(803, 15).Value = 1
$ws.Cells.Item(803, 16).Formula = "=IF(M803>I803,1,0)"
$ws.Cells.Item(803, 17).Formula = "=IF(P803=1,(M803-I803), """")"
$ws.Cells.Item(803, 18).Formula = "=IF(M803<I803, 1, 0)"
$ws.Cells.Item(803, 19).Formula = "=IF(R803=1,(I803-M803),"""")"
$ws.Cells.Item(803, 20).Formula = "=IF(M803=I803,1,0)"

# Row 804: TOR @ CLE
$ws.Cells.Item(804, 1).Value = 45700
$ws.Cells.Item(804, 2).Value = "Wednesday"
$ws.Cells.Item(804, 3).Value = "CLE"
$ws.Cells.Item(804, 4).Value = "TOR"
$ws.Cells.Item(804, 5).Value = 0
$ws.Cells.Item(804, 6).Value = "Marc Davis"
$ws.Cells.Item(804, 7).Value = "Natalie Sago"
$ws.Cells.Item(804, 8).Value = "Brandon Schwab"
$ws.Cells.Item(804, 9).Value = 235.5
$ws.Cells.Item(804, 10).Value = 12.5
$ws.Cells.Item(804, 11).Value = 131
$ws.Cells.Item(804, 12).Value = 108
$ws.Cells.Item(804, 13).Formula = "=K804+L804"
$ws.Cells.Item(804, 14).Formula = "=(L804-K804)*-1"
$ws.Cells.Item(804, 15).Value = 1
$ws.Cells.Item(804, 16).Formula = "=IF(M804>I804,1,0)"
$ws.Cells.Item(804, 17).Formula = "=IF(P804=1,(M804-I804), """")"
$ws.Cells.Item(804, 18).Formula = "=IF(M804<I804, 1, 0)"
$ws.Cells.Item(804, 19).Formula = "=IF(R804=1,(I804-M804),"""")"
$ws.Cells.Item(804, 20).Formula = "=IF(M804=I804,1,0)"

# Row 805: NYK @ ATL
$ws.Cells.Item(805, 1).Value = 45700
$ws.Cells.Item(805, 2).Value = "Wednesday"
$ws.Cells.Item(805, 3).Value = "ATL"
$ws.Cells.Item(805, 4).Value = "NYK"
$ws.Cells.Item(805, 5).Value = 0
$ws.Cells.Item(805, 6).Value = "Mark Lindsay"
$ws.Cells.Item(805, 7).Value = "Nate Green"
$ws.Cells.Item(805, 8).Value = "Michael Smith"
$ws.Cells.Item(805, 9).Value = 240
$ws.Cells.Item(805, 10).Value = -8.5
$ws.Cells.Item(805, 11).Value = 148
$ws.Cells.Item(805, 12).Value = 149
$ws.Cells.Item(805, 13).Formula = "=K805+L805"
$ws.Cells.Item(805, 14).Formula = "=(L805-K805)*-1"
$ws.Cells.Item(805, 15).Value = 1
$ws.Cells.Item(805, 16).Formula = "=IF(M805>I805,1,0)"
$ws.Cells.Item(805, 17).Formula = "=IF(P805=1,(M805-I805), """")"
$ws.Cells.Item(805, 18).Formula = "=IF(M805<I805, 1, 0)"
$ws.Cells.Item(805, 19).Formula = "=IF(R805=1,(I805-M805),"""")"
$ws.Cells.Item(805, 20).Formula = "=IF(M805=I805,1,0)"

# Row 806: NOP @ SAC
$ws.Cells.Item(806, 1).Value = 45700
$ws.Cells.Item(806, 2).Value = "Wednesday"
$ws.Cells.Item(806, 3).Value = "SAC"
$ws.Cells.Item(806, 4).Value = "NOP"
$ws.Cells.Item(806, 5).Value = 0
$ws.Cells.Item(806, 6).Value = "Gediminas Petraitis"
$ws.Cells.Item(806, 7).Value = "Evan Scott"
$ws.Cells.Item(806, 8).Value = "Matt Myers"
$ws.Cells.Item(806, 9).Value = 238.5
$ws.Cells.Item(806, 10).Value = 4
$ws.Cells.Item(806, 11).Value = 119
$ws.Cells.Item(806, 12).Value = 111
$ws.Cells.Item(806, 13).Formula = "=K806+L806"
$ws.Cells.Item(806, 14).Formula = "=(L806-K806)*-1"
$ws.Cells.Item(806, 15).Value = 1
$ws.Cells.Item(806, 16).Formula = "=IF(M806>I806,1,0)"
$ws.Cells.Item(806, 17).Formula = "=IF(P806=1,(M806-I806), """")"
$ws.Cells.Item(806, 18).Formula = "=IF(M806<I806, 1, 0)"
$ws.Cells.Item(806, 19).Formula = "=IF(R806=1,(I806-M806),"""")"
$ws.Cells.Item(806, 20).Formula = "=IF(M806=I806,1,0)"

# Row 807: OKC @ MIA
$ws.Cells.Item(807, 1).Value = 45700
$ws.Cells.Item(807, 2).Value = "Wednesday"
$ws.Cells.Item(807, 3).Value = "MIA"
$ws.Cells.Item(807, 4).Value = "OKC"
$ws.Cells.Item(807, 5).Value = 0
$ws.Cells.Item(807, 6).Value = "Scott Foster"
$ws.Cells.Item(807, 7).Value = "Karl Lane"
$ws.Cells.Item(807, 8).Value = "Brett Nansel"
$ws.Cells.Item(807, 9).Value = 219.5
$ws.Cells.Item(807, 10).Value = -13.5
$ws.Cells.Item(807, 11).Value = 101
$ws.Cells.Item(807, 12).Value = 115
$ws.Cells.Item(807, 13).Formula = "=K807+L807"
$ws.Cells.Item(807, 14).Formula = "=(L807-K807)*-1"
$ws.Cells.Item(807, 15).Value = 1
$ws.Cells.Item(807, 16).Formula = "=IF(M807>I807,1,0)"
$ws.Cells.Item(807, 17).Formula = "=IF(P807=1,(M807-I807), """")"
$ws.Cells.Item(807, 18).Formula = "=IF(M807<I807, 1, 0)"
$ws.Cells.Item(807, 19).Formula = "=IF(R807=1,(I807-M807),"""")"
$ws.Cells.Item(807, 20).Formula = "=IF(M807=I807,1,0)"

# Row 808: MIN @ MIL
$ws.Cells.Item(808, 1).Value = 45700
$ws.Cells.Item(808, 2).Value = "Wednesday"
$ws.Cells.Item(808, 3).Value = "MIL"
$ws.Cells.Item(808, 4).Value = "MIN"
$ws.Cells.Item(808, 5).Value = 0
$ws.Cells.Item(808, 6).Value = "Mitchell Ervin"
$ws.Cells.Item(808, 7).Value = "Marat Kogut"
$ws.Cells.Item(808, 8).Value = "Tom Washington"
$ws.Cells.Item(808, 9).Value = 226
$ws.Cells.Item(808, 10).Value = -5.5
$ws.Cells.Item(808, 11).Value = 103
$ws.Cells.Item(808, 12).Value = 101
$ws.Cells.Item(808, 13).Formula = "=K808+L808"
$ws.Cells.Item(808, 14).Formula = "=(L808-K808)*-1"
$ws.Cells.Item(808, 15).Value = 1
$ws.Cells.Item(808, 16).Formula = "=IF(M808>I808,1,0)"
$ws.Cells.Item(808, 17).Formula = "=IF(P808=1,(M808-I808), """")"
$ws.Cells.Item(808, 18).Formula = "=IF(M808<I808, 1, 0)"
$ws.Cells.Item(808, 19).Formula = "=IF(R808=1,(I808-M808),"""")"
$ws.Cells.Item(808, 20).Formula = "=IF(M808=I808,1,0)"

# Row 809: CHI @ DET
$ws.Cells.Item(809, 1).Value = 45700
$ws.Cells.Item(809, 2).Value = "Wednesday"
$ws.Cells.Item(809, 3).Value = "DET"
$ws.Cells.Item(809, 4).Value = "CHI"
$ws.Cells.Item(809, 5).Value = 0
$ws.Cells.Item(809, 6).Value = "Pat Fraher"
$ws.Cells.Item(809, 7).Value = "Jason Goldenberg"
$ws.Cells.Item(809, 8).Value = "Brandon Adair"
$ws.Cells.Item(809, 9).Value = 236.5
$ws.Cells.Item(809, 10).Value = 3.5
$ws.Cells.Item(809, 11).Value = 128
$ws.Cells.Item(809, 12).Value = 110
$ws.Cells.Item(809, 13).Formula = "=K809+L809"
$ws.Cells.Item(809, 14).Formula = "=(L809-K809)*-1"
$ws.Cells.Item(809, 15).Value = 1
$ws.Cells.Item(809, 16).Formula = "=IF(M809>I809,1,0)"
$ws.Cells.Item(809, 17).Formula = "=IF(P809=1,(M809-I809), """")"
$ws.Cells.Item(809, 18).Formula = "=IF(M809<I809, 1, 0)"
$ws.Cells.Item(809, 19).Formula = "=IF(R809=1,(I809-M809),"""")"
$ws.Cells.Item(809, 20).Formula = "=IF(M809=I809,1,0)"

# Row 810: HOU @ PHX
$ws.Cells.Item(810, 1).Value = 45700
$ws.Cells.Item(810, 2).Value = "Wednesday"
$ws.Cells.Item(810, 3).Value = "PHX"
$ws.Cells.Item(810, 4).Value = "HOU"
$ws.Cells.Item(810, 5).Value = 0
$ws.Cells.Item(810, 6).Value = "Zach Zarba"
$ws.Cells.Item(810, 7).Value = "John Butler"
$ws.Cells.Item(810, 8).Value = "Simone Jelks"
$ws.Cells.Item(810, 9).Value = 223.5
$ws.Cells.Item(810, 10).Value = -6
$ws.Cells.Item(810, 11).Value = 111
$ws.Cells.Item(810, 12).Value = 119
$ws.Cells.Item(810, 13).Formula = "=K810+L810"
$ws.Cells.Item(810, 14).Formula = "=(L810-K810)*-1"
$ws.Cells.Item(810, 15).Value = 1
$ws.Cells.Item(810, 16).Formula = "=IF(M810>I810,1,0)"
$ws.Cells.Item(810, 17).Formula = "=IF(P810=1,(M810-I810), """")"
$ws.Cells.Item(810, 18).Formula = "=IF(M810<I810, 1, 0)"
$ws.Cells.Item(810, 19).Formula = "=IF(R810=1,(I810-M810),"""")"
$ws.Cells.Item(810, 20).Formula = "=IF(M810=I810,1,0)"

# Row 811: DEN @ POR
$ws.Cells.Item(811, 1).Value = 45700
$ws.Cells.Item(811, 2).Value = "Wednesday"
$ws.Cells.Item(811, 3).Value = "POR"
$ws.Cells.Item(811, 4).Value = "DEN"
$ws.Cells.Item(811, 5).Value = 0
$ws.Cells.Item(811, 6).Value = "Tyler Ford"
$ws.Cells.Item(811, 7).Value = "Andy Nagy"
$ws.Cells.Item(811, 8).Value = "CJ Washington"
$ws.Cells.Item(811, 9).Value = 231.5
$ws.Cells.Item(811, 10).Value = -11.5
$ws.Cells.Item(811, 11).Value = 121
$ws.Cells.Item(811, 12).Value = 132
$ws.Cells.Item(811, 13).Formula = "=K811+L811"
$ws.Cells.Item(811, 14).Formula = "=(L811-K811)*-1"
$ws.Cells.Item(811, 15).Value = 1
$ws.Cells.Item(811, 16).Formula = "=IF(M811>I811,1,0)"
$ws.Cells.Item(811, 17).Formula = "=IF(P811=1,(M811-I811), """")"
$ws.Cells.Item(811, 18).Formula = "=IF(M811<I811, 1, 0)"
$ws.Cells.Item(811, 19).Formula = "=IF(R811=1,(I811-M811),"""")"
$ws.Cells.Item(811, 20).Formula = "=IF(M811=I811,1,0)"

# Row 812: UTA @ LAL
$ws.Cells.Item(812, 1).Value = 45700
$ws.Cells.Item(812, 2).Value = "Wednesday"
$ws.Cells.Item(812, 3).Value = "LAL"
$ws.Cells.Item(812, 4).Value = "UTA"
$ws.Cells.Item(812, 5).Value = 0
$ws.Cells.Item(812, 6).Value = "Tre Maddox"
$ws.Cells.Item(812, 7).Value = "JT Orr"
$ws.Cells.Item(812, 8).Value = "Mousa Dagher"
$ws.Cells.Item(812, 9).Value = 237.5
$ws.Cells.Item(812, 10).Value = 7.5
$ws.Cells.Item(812, 11).Value = 119
$ws.Cells.Item(812, 12).Value = 131
$ws.Cells.Item(812, 13).Formula = "=K812+L812"
$ws.Cells.Item(812, 14).Formula = "=(L812-K812)*-1"
$ws.Cells.Item(812, 15).Value = 1
$ws.Cells.Item(812, 16).Formula = "=IF(M812>I812,1,0)"
$ws.Cells.Item(812, 17).Formula = "=IF(P812=1,(M812-I812), """")"
$ws.Cells.Item(812, 18).Formula = "=IF(M812<I812, 1, 0)"
$ws.Cells.Item(812, 19).Formula = "=IF(R812=1,(I812-M812),"""")"
$ws.Cells.Item(812, 20).Formula = "=IF(M812=I812,1,0)"

# Row 813: DAL @ GSW
$ws.Cells.Item(813, 1).Value = 45700
$ws.Cells.Item(813, 2).Value = "Wednesday"
$ws.Cells.Item(813, 3).Value = "GSW"
$ws.Cells.Item(813, 4).Value = "DAL"
$ws.Cells.Item(813, 5).Value = 0
$ws.Cells.Item(813, 6).Value = "James Williams"
$ws.Cells.Item(813, 7).Value = "Jacyn Goble"
$ws.Cells.Item(813, 8).Value = "Jonathan Sterling"
$ws.Cells.Item(813, 9).Value = 232.5
$ws.Cells.Item(813, 10).Value = 3.5
$ws.Cells.Item(813, 11).Value = 107
$ws.Cells.Item(813, 12).Value = 111
$ws.Cells.Item(813, 13).Formula = "=K813+L813"
$ws.Cells.Item(813, 14).Formula = "=(L813-K813)*-1"
$ws.Cells.Item(813, 15).Value = 1
$ws.Cells.Item(813, 16).Formula = "=IF(M813>I813,1,0)"
$ws.Cells.Item(813, 17).Formula = "=IF(P813=1,(M813-I813), """")"
$ws.Cells.Item(813, 18).Formula = "=IF(M813<I813, 1, 0)"
$ws.Cells.Item(813, 19).Formula = "=IF(R813=1,(I813-M813),"""")"
$ws.Cells.Item(813, 20).Formula = "=IF(M813=I813,1,0)"

# Row 814: LAC @ MEM
$ws.Cells.Item(814, 1).Value = 45700
$ws.Cells.Item(814, 2).Value = "Wednesday"
$ws.Cells.Item(814, 3).Value = "MEM"
$ws.Cells.Item(814, 4).Value = "LAC"
$ws.Cells.Item(814, 5).Value = 0
$ws.Cells.Item(814, 6).Value = "Kevin Cutler"
$ws.Cells.Item(814, 7).Value = "Scott Twardoski"
$ws.Cells.Item(814, 8).Value = "Derek Richardson"
$ws.Cells.Item(814, 9).Value = 235.5
$ws.Cells.Item(814, 10).Value = -3.5
$ws.Cells.Item(814, 11).Value = 114
$ws.Cells.Item(814, 12).Value = 128
$ws.Cells.Item(814, 13).Formula = "=K814+L814"
$ws.Cells.Item(814, 14).Formula = "=(L814-K814)*-1"
$ws.Cells.Item(814, 15).Value = 1
$ws.Cells.Item(814, 16).Formula = "=IF(M814>I814,1,0)"
$ws.Cells.Item(814, 17).Formula = "=IF(P814=1,(M814-I814), """")"
$ws.Cells.Item(814, 18).Formula = "=IF(M814<I814, 1, 0)"
$ws.Cells.Item(814, 19).Formula = "=IF(R814=1,(I814-M814),"""")"
$ws.Cells.Item(814, 20).Formula = "=IF(M814=I814,1,0)"

# Match the final selection state left by the author after entering the data
$ws.Range("V803").Select()

Write-Host "Added game data for February 12, 2025 (rows 800-814)"